# Apply cryptos list update per commit "Updated cryptos list on Sat Apr 13 15:44:13 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.751.54"
$ws.Range("E2").Value = "'  -2.47%  "

$ws.Range("D3").Value = "'3.280.42"

$ws.Range("E4").Value = "'  +0.03%  "

$ws.Range("D5").Value = "'595.11"
$ws.Range("E5").Value = "'  -2.68%  "

$ws.Range("D6").Value = "'151.27"
$ws.Range("E6").Value = "'  -9.78%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.17%  "

$ws.Range("D8").Value = "'3.276.00"
$ws.Range("E8").Value = "'  -4.65%  "

$ws.Range("D9").Value = "'0.546"
$ws.Range("E9").Value = "'  -8.20%  "

$ws.Range("E10").Value = "'  -10.80%  "

$ws.Range("D11").Value = "'6.69"
$ws.Range("E11").Value = "'  -5.02%  "

$ws.Range("D12").Value = "'0.506"
$ws.Range("E12").Value = "'  -10.21%  "

$ws.Range("D13").Value = "'38.69"
$ws.Range("E13").Value = "'  -12.69%  "

$ws.Range("D14").Value = "'0.0000247"
$ws.Range("E14").Value = "'  -8.23%  "

$ws.Range("D15").Value = "'3.806.44"
$ws.Range("E15").Value = "'  -4.89%  "

$ws.Range("D16").Value = "'67.753.65"
$ws.Range("E16").Value = "'  -2.60%  "

$ws.Range("D17").Value = "'3.283.32"
$ws.Range("E17").Value = "'  -4.77%  "

$ws.Range("D18").Value = "'532.82"
$ws.Range("E18").Value = "'  -8.46%  "

$ws.Range("E19").Value = "'  -5.57%  "

$ws.Range("D20").Value = "'7.17"
$ws.Range("E20").Value = "'  -12.05%  "

$ws.Range("D21").Value = "'15.04"
$ws.Range("E21").Value = "'  -12.41%  "

$ws.Range("D22").Value = "'0.759"
$ws.Range("E22").Value = "'  -10.41%  "

$ws.Range("D23").Value = "'7.86"
$ws.Range("E23").Value = "'  -11.61%  "

$ws.Range("D24").Value = "'85.94"
$ws.Range("E24").Value = "'  -10.22%  "

$ws.Range("D25").Value = "'13.58"
$ws.Range("E25").Value = "'  -10.54%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "'  -0.18%  "

$ws.Range("E27").Value = "'  -10.63%  "

$ws.Range("D28").Value = "'8.12"
$ws.Range("E28").Value = "'  -5.93%  "

$ws.Range("E29").Value = "'  -10.78%  "

$ws.Range("D30").Value = "'29.18"
$ws.Range("E30").Value = "'  -10.98%  "

$ws.Range("E31").Value = "'  -2.77%  "

$ws.Range("D32").Value = "'2.68"
$ws.Range("E32").Value = "'  -4.38%  "

$ws.Range("D33").Value = "'6.64"
$ws.Range("E33").Value = "'  -15.03%  "

$ws.Range("E34").Value = "'  -12.15%  "

$ws.Range("B35").Value = "'OKB"
$ws.Range("C35").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'57.64"
$ws.Range("E35").Value = "'  +2.81%  "

$ws.Range("B36").Value = "'Bittensor"
$ws.Range("C36").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "'516.54"
$ws.Range("E36").Value = "'  -10.26%  "

$ws.Range("E37").Value = "'  -0.05%  "

$ws.Range("D38").Value = "'0.0450"
$ws.Range("E38").Value = "'  -5.83%  "

$ws.Range("D39").Value = "'0.0859"
$ws.Range("E39").Value = "'  -10.07%  "

$ws.Range("D40").Value = "'9.01"
$ws.Range("E40").Value = "'  -14.38%  "

$ws.Range("E41").Value = "'  -9.62%  "

$ws.Range("D42").Value = "'2.80"
$ws.Range("E42").Value = "'  -10.96%  "

$ws.Range("D43").Value = "'2.953.87"
$ws.Range("E43").Value = "'  -8.84%  "

$ws.Range("E44").Value = "'  -8.66%  "

$ws.Range("D45").Value = "'0.0₃0589"
$ws.Range("E45").Value = "'  -14.01%  "

$ws.Range("E46").Value = "'  -6.78%  "

$ws.Range("D47").Value = "'26.80"
$ws.Range("E47").Value = "'  -13.93%  "

$ws.Range("E48").Value = "'  -0.11%  "

$ws.Range("E49").Value = "'  -15.21%  "

$ws.Range("E50").Value = "'  -9.85%  "

$ws.Range("D51").Value = "'123.80"
$ws.Range("E51").Value = "'  -7.66%  "

Write-Host "Applied crypto price/volume updates"
